# Adds the "other candidates at 23C" data block (rows 38-59) to Sheet1.
# Two experimental conditions (LB_CAM_KAN, then LB_RIF_CAM_KAN) are recorded
# for 11 new constructs, mirroring the pre-existing layout used for rows 2-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# New constructs tested at 23C, in the order they appear in the sheet.
$constructs = @("CDA1", "A2", "CDAL7", "PyrD", "GsDA", "CDAL4", "TAD3", "CDAL1", "CDAL2", "CDAL5", "CDAL6")

# L-column (raw count) values for the LB_CAM_KAN condition rows 38-48.
$countsCamKan = @(420, 295, 305, 341, 287, 449, 327, 297, 372, 404, 239)

# L-column (raw count) values for the LB_RIF_CAM_KAN condition rows 49-59.
$countsRifCamKan = @(382, 380, 344, 504, 346, 482, 509, 485, 840, 609, 468)

$n = $constructs.Length
$firstRow = 38
$secondRow = 49

for ($i = 0; $i -lt $n; $i++) {
    $rowCam = $firstRow + $i
    $rowRif = $secondRow + $i

    # --- LB_CAM_KAN row: reuse formatting from row 32 (same condition pattern) ---
    $ws.Range("B32").Copy()
    $ws.Range("B$rowCam").PasteSpecial($xlPasteFormats)
    $ws.Range("C32").Copy()
    $ws.Range("C$rowCam").PasteSpecial($xlPasteFormats)
    $ws.Range("E32").Copy()
    $ws.Range("E$rowCam").PasteSpecial($xlPasteFormats)

    $ws.Range("A$rowCam").Value = "LZ"
    $ws.Range("B$rowCam").Value = 45370
    $ws.Range("C$rowCam").Value = $constructs[$i]
    $ws.Range("D$rowCam").Value = "bl21_de3_delta_ung"
    $ws.Range("E$rowCam").Value = "LB_CAM_KAN"
    $ws.Range("G$rowCam").Value = 23
    $ws.Range("H$rowCam").Value = 1
    $ws.Range("I$rowCam").Value = 1
    $ws.Range("K$rowCam").Formula = "=10^-6*0.05"
    $ws.Range("L$rowCam").Value = $countsCamKan[$i]
    $ws.Range("M$rowCam").Formula = "=L$rowCam/K$rowCam"

    # --- LB_RIF_CAM_KAN row: reuse formatting from row 33 (same condition pattern) ---
    $ws.Range("B33").Copy()
    $ws.Range("B$rowRif").PasteSpecial($xlPasteFormats)
    $ws.Range("C33").Copy()
    $ws.Range("C$rowRif").PasteSpecial($xlPasteFormats)
    $ws.Range("E33").Copy()
    $ws.Range("E$rowRif").PasteSpecial($xlPasteFormats)

    $ws.Range("A$rowRif").Value = "LZ"
    $ws.Range("B$rowRif").Value = 45370
    $ws.Range("C$rowRif").Value = $constructs[$i]
    $ws.Range("D$rowRif").Value = "bl21_de3_delta_ung"
    $ws.Range("E$rowRif").Value = "LB_RIF_CAM_KAN"
    $ws.Range("G$rowRif").Value = 23
    $ws.Range("H$rowRif").Value = 1
    $ws.Range("I$rowRif").Value = 1
    $ws.Range("K$rowRif").Value = 1
    $ws.Range("L$rowRif").Value = $countsRifCamKan[$i]
    $ws.Range("M$rowRif").Formula = "=L$rowRif/K$rowRif"
    $ws.Range("N$rowRif").Formula = "=M$rowRif/M$rowCam"
}

# --- Sync the sheet view to match where the author left the cursor ---
$win = $excel.ActiveWindow
$win.Zoom = 57
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("K56").Select()
